$wb = $excel.ActiveWorkbook

# Updated market-price / profit figures pulled by the scheduled runner.
# Each block targets one worksheet row; values come from the latest
# Universalis price snapshot, recomputed profit columns follow.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 145.8
$ws.Range("I12").Value = 106.333336
$ws.Range("J12").Value = 205
$ws.Range("K12").Value = 106.333336
$ws.Range("L12").Value = 205
$ws.Range("M12").Value = 63.666664
$ws.Range("N12").Value = -545

$ws.Range("H40").Value = 2099.9333
$ws.Range("I40").Value = 2138.3845
$ws.Range("J40").Value = 1850
$ws.Range("K40").Value = 2138.3845
$ws.Range("L40").Value = 1850
$ws.Range("M40").Value = -1963.3845
$ws.Range("N40").Value = -2200

$ws.Range("H55").Value = 4273.815
$ws.Range("I55").Value = 1245.7778
$ws.Range("J55").Value = 5787.8335
$ws.Range("K55").Value = 1245.7778
$ws.Range("L55").Value = 5787.8335
$ws.Range("M55").Value = -1031.7778
$ws.Range("N55").Value = -6215.8335

$ws.Range("H64").Value = 4072.65
$ws.Range("I64").Value = 3428.5715
$ws.Range("J64").Value = 5575.5
$ws.Range("K64").Value = 3428.5715
$ws.Range("L64").Value = 5575.5
$ws.Range("M64").Value = -3180.5715
$ws.Range("N64").Value = -6071.5

$ws.Range("H67").Value = 4072.65
$ws.Range("I67").Value = 3428.5715
$ws.Range("J67").Value = 5575.5
$ws.Range("K67").Value = 3428.5715
$ws.Range("L67").Value = 5575.5
$ws.Range("M67").Value = -2570.5715
$ws.Range("N67").Value = -7291.5

$ws.Range("H116").Value = 2984.5
$ws.Range("I116").Value = 2866.6667
$ws.Range("J116").Value = 3055.2
$ws.Range("K116").Value = 2866.6667
$ws.Range("L116").Value = 3055.2
$ws.Range("M116").Value = 575.3332999999998
$ws.Range("N116").Value = -9939.200000000001

$ws.Range("H138").Value = 1909.4546
$ws.Range("I138").Value = 1819.4
$ws.Range("K138").Value = 5458.200000000001
$ws.Range("M138").Value = -318.2000000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 18667
$ws.Range("I22").Value = 3000.5
$ws.Range("K22").Value = 3000.5
$ws.Range("M22").Value = -2701.5

$ws.Range("H32").Value = 4226.4307
$ws.Range("I32").Value = 2698.3076
$ws.Range("K32").Value = 2698.3076
$ws.Range("M32").Value = -2411.3076

$ws.Range("H63").Value = 9092.857
$ws.Range("I63").Value = 9561.538
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 9561.538
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -8875.538
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 9092.857
$ws.Range("I66").Value = 9561.538
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 47807.69
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -44375.69
$ws.Range("N66").Value = -21864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 221.88889
$ws.Range("I64").Value = 148.5
$ws.Range("J64").Value = 280.6
$ws.Range("K64").Value = 148.5
$ws.Range("L64").Value = 280.6
$ws.Range("M64").Value = 76.5
$ws.Range("N64").Value = -730.6

$ws.Range("H67").Value = 221.88889
$ws.Range("I67").Value = 148.5
$ws.Range("J67").Value = 280.6
$ws.Range("K67").Value = 148.5
$ws.Range("L67").Value = 280.6
$ws.Range("M67").Value = 631.5
$ws.Range("N67").Value = -1840.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 70005
$ws.Range("J23").Value = 70005
$ws.Range("L23").Value = 70005
$ws.Range("N23").Value = -70485

$ws.Range("H27").Value = 70005
$ws.Range("J27").Value = 70005
$ws.Range("L27").Value = 70005
$ws.Range("N27").Value = -70389

$ws.Range("H31").Value = 2163.1875
$ws.Range("I31").Value = 1899.7858
$ws.Range("J31").Value = 4007
$ws.Range("K31").Value = 1899.7858
$ws.Range("L31").Value = 4007
$ws.Range("M31").Value = -1604.7858
$ws.Range("N31").Value = -4597

$ws.Range("H34").Value = 2163.1875
$ws.Range("I34").Value = 1899.7858
$ws.Range("J34").Value = 4007
$ws.Range("K34").Value = 1899.7858
$ws.Range("L34").Value = 4007
$ws.Range("M34").Value = -1697.7858
$ws.Range("N34").Value = -4411

$ws.Range("H58").Value = 2364.4666
$ws.Range("I58").Value = 1450.4
$ws.Range("J58").Value = 3278.5334
$ws.Range("K58").Value = 1450.4
$ws.Range("L58").Value = 3278.5334
$ws.Range("M58").Value = -1247.4
$ws.Range("N58").Value = -3684.5334

$ws.Range("H136").Value = 2364.4666
$ws.Range("I136").Value = 1450.4
$ws.Range("J136").Value = 3278.5334
$ws.Range("K136").Value = 4351.200000000001
$ws.Range("L136").Value = 9835.600199999999
$ws.Range("M136").Value = -1801.200000000001
$ws.Range("N136").Value = -14935.6002

$ws.Range("H137").Value = 33254
$ws.Range("J137").Value = 33254
$ws.Range("L137").Value = 33254
$ws.Range("N137").Value = -43454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1398.1765
$ws.Range("I34").Value = 129.8
$ws.Range("J34").Value = 1926.6666
$ws.Range("K34").Value = 389.4
$ws.Range("L34").Value = 5779.9998
$ws.Range("M34").Value = -305.4
$ws.Range("N34").Value = -5947.9998

$ws.Range("H39").Value = 6536.2974
$ws.Range("J39").Value = 4441.9355
$ws.Range("L39").Value = 13325.8065
$ws.Range("N39").Value = -13913.8065

$ws.Range("H104").Value = 3553.9092
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 3553.9092
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 10661.7276
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -15903.7276

$ws.Range("H117").Value = 1980
$ws.Range("J117").Value = 2966.6667
$ws.Range("L117").Value = 8900.000100000001
$ws.Range("N117").Value = -15784.0001

$ws.Range("H131").Value = 2448.2666
$ws.Range("I131").Value = 450
$ws.Range("J131").Value = 2560.845
$ws.Range("K131").Value = 1350
$ws.Range("L131").Value = 7682.535
$ws.Range("M131").Value = 3690
$ws.Range("N131").Value = -17762.535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H24").Value = 126892.305
$ws.Range("J24").Value = 41600
$ws.Range("L24").Value = 41600
$ws.Range("N24").Value = -41946

$ws.Range("H39").Value = 15600
$ws.Range("J39").Value = 15400
$ws.Range("L39").Value = 15400
$ws.Range("N39").Value = -16464

$ws.Range("H80").Value = 1540.7059
$ws.Range("I80").Value = 1648.7142
$ws.Range("J80").Value = 1465.1
$ws.Range("K80").Value = 1648.7142
$ws.Range("L80").Value = 1465.1
$ws.Range("M80").Value = -650.7141999999999
$ws.Range("N80").Value = -3461.1

$ws.Range("H83").Value = 1540.7059
$ws.Range("I83").Value = 1648.7142
$ws.Range("J83").Value = 1465.1
$ws.Range("K83").Value = 8243.571
$ws.Range("L83").Value = 7325.5
$ws.Range("M83").Value = -3251.571
$ws.Range("N83").Value = -17309.5

$ws.Range("H132").Value = 5701.9165
$ws.Range("I132").Value = 3944.4285
$ws.Range("J132").Value = 8162.4
$ws.Range("K132").Value = 11833.2855
$ws.Range("L132").Value = 24487.2
$ws.Range("M132").Value = -9303.2855
$ws.Range("N132").Value = -29547.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45309.914
$ws.Range("I7").Value = 72741.71000000001
$ws.Range("J7").Value = 2638.2222
$ws.Range("K7").Value = 72741.71000000001
$ws.Range("L7").Value = 2638.2222
$ws.Range("M7").Value = -72629.71000000001
$ws.Range("N7").Value = -2862.2222

$ws.Range("H22").Value = 1521.8572
$ws.Range("I22").Value = 661
$ws.Range("J22").Value = 1866.2
$ws.Range("K22").Value = 661
$ws.Range("L22").Value = 1866.2
$ws.Range("M22").Value = -366
$ws.Range("N22").Value = -2456.2

$ws.Range("H27").Value = 1521.8572
$ws.Range("I27").Value = 661
$ws.Range("J27").Value = 1866.2
$ws.Range("K27").Value = 661
$ws.Range("L27").Value = 1866.2
$ws.Range("M27").Value = -554
$ws.Range("N27").Value = -2080.2

$ws.Range("H46").Value = 201300.2
$ws.Range("I46").Value = 251000.25
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 251000.25
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = -250812.25
$ws.Range("N46").Value = -2876

$ws.Range("H68").Value = 2502.9
$ws.Range("J68").Value = 2754.8333
$ws.Range("L68").Value = 2754.8333
$ws.Range("N68").Value = -4252.8333

$ws.Range("H71").Value = 2502.9
$ws.Range("J71").Value = 2754.8333
$ws.Range("L71").Value = 13774.1665
$ws.Range("N71").Value = -21262.1665

$ws.Range("H82").Value = 1253.0625
$ws.Range("I82").Value = 1099.8572
$ws.Range("J82").Value = 1372.2222
$ws.Range("K82").Value = 1099.8572
$ws.Range("L82").Value = 1372.2222
$ws.Range("M82").Value = -738.8571999999999
$ws.Range("N82").Value = -2094.2222

$ws.Range("H85").Value = 1253.0625
$ws.Range("I85").Value = 1099.8572
$ws.Range("J85").Value = 1372.2222
$ws.Range("K85").Value = 1099.8572
$ws.Range("L85").Value = 1372.2222
$ws.Range("M85").Value = 148.1428000000001
$ws.Range("N85").Value = -3868.2222

$ws.Range("H126").Value = 45309.914
$ws.Range("I126").Value = 72741.71000000001
$ws.Range("J126").Value = 2638.2222
$ws.Range("K126").Value = 218225.13
$ws.Range("L126").Value = 7914.6666
$ws.Range("M126").Value = -215755.13
$ws.Range("N126").Value = -12854.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 37509.25
$ws.Range("I31").Value = 4999.5
$ws.Range("J31").Value = 70019
$ws.Range("K31").Value = 4999.5
$ws.Range("L31").Value = 70019
$ws.Range("M31").Value = -4651.5
$ws.Range("N31").Value = -70715

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2208

$ws.Range("H122").Value = 57382.777
$ws.Range("I122").Value = 57382.777
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 172148.331
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -169698.331
$ws.Range("N122").ClearContents()
